$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet gets a new (blank) column inserted
# immediately before the old "Late" column (column N), pushing the
# existing Late / Outstanding(heading) / Outstanding columns one to the
# right (N->O, O->P, P->Q).
$ws = $wb.Worksheets.Item("Repayment schedule")

$null = $ws.Columns("N").Insert()

# Give the freshly inserted column the same effective width as column M
# (a manually-set width, not an auto "best fit" one).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with the given cell
# selected - this also clears the previous tab selection (NewLoanInput).
$ws.Activate()
$null = $ws.Range("K16").Select()
